# Auto-generated edit script applying scheduled market-price/profit refresh
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# ALC row 40
$ws1.Range("H40").Value = 2417.8
$ws1.Range("I40").Value = 2099.5
$ws1.Range("J40").Value = 2630
$ws1.Range("K40").Value = 2099.5
$ws1.Range("L40").Value = 2630
$ws1.Range("M40").Value = -1924.5
$ws1.Range("N40").Value = -2980

# ALC row 106
$ws1.Range("H106").Value = 33335840
$ws1.Range("I106").Value = 41669320
$ws1.Range("K106").Value = 41669320
$ws1.Range("M106").Value = -41668689

# ALC row 132
$ws1.Range("H132").Value = 2426.8928
$ws1.Range("I132").Value = 1876.8422
$ws1.Range("K132").Value = 5630.5266
$ws1.Range("M132").Value = -3100.5266

# ALC row 137
$ws1.Range("H137").Value = 3224.2156
$ws1.Range("I137").Value = 1207.5
$ws1.Range("K137").Value = 3622.5
$ws1.Range("M137").Value = -1072.5

# ALC row 138
$ws1.Range("H138").Value = 3502.859
$ws1.Range("I138").Value = 1055.0476
$ws1.Range("K138").Value = 3165.142800000001
$ws1.Range("M138").Value = 1974.857199999999

# ARM row 2
$ws2.Range("H2").Value = 25643824
$ws2.Range("I2").Value = 32260618
$ws2.Range("K2").Value = 32260618
$ws2.Range("M2").Value = -32260505

# ARM row 19
$ws2.Range("H19").Value = 7755.8335
$ws2.Range("J19").Value = 1009
$ws2.Range("L19").Value = 1009
$ws2.Range("N19").Value = -1467

# ARM row 61
$ws2.Range("H61").Value = 3311.7568
$ws2.Range("I61").Value = 2993.4517
$ws2.Range("J61").Value = 4956.3335
$ws2.Range("K61").Value = 2993.4517
$ws2.Range("L61").Value = 4956.3335
$ws2.Range("M61").Value = -2781.4517
$ws2.Range("N61").Value = -5380.3335

# ARM row 102
$ws2.Range("H102").Value = 9278.933999999999
$ws2.Range("I102").Value = 9185.77
$ws2.Range("K102").Value = 9185.77
$ws2.Range("M102").Value = -7563.77

# ARM row 116
$ws2.Range("H116").Value = 25643824
$ws2.Range("I116").Value = 32260618
$ws2.Range("K116").Value = 32260618
$ws2.Range("M116").Value = -32258324

# ARM row 132
$ws2.Range("H132").Value = 3922.9614
$ws2.Range("I132").Value = 3399.2
$ws2.Range("J132").Value = 5668.8335
$ws2.Range("K132").Value = 10197.6
$ws2.Range("L132").Value = 17006.5005
$ws2.Range("M132").Value = -7667.599999999999
$ws2.Range("N132").Value = -22066.5005

# ARM row 136
$ws2.Range("H136").Value = 3311.7568
$ws2.Range("I136").Value = 2993.4517
$ws2.Range("J136").Value = 4956.3335
$ws2.Range("K136").Value = 8980.355100000001
$ws2.Range("L136").Value = 14869.0005
$ws2.Range("M136").Value = -6430.355100000001
$ws2.Range("N136").Value = -19969.0005

# BSM row 3
$ws3.Range("H3").Value = 25643824
$ws3.Range("I3").Value = 32260618
$ws3.Range("K3").Value = 32260618
$ws3.Range("M3").Value = -32260504

# BSM row 22
$ws3.Range("H22").Value = 406550.25
$ws3.Range("I22").Value = 767.5
$ws3.Range("K22").Value = 767.5
$ws3.Range("M22").Value = -594.5

# BSM row 99
$ws3.Range("H99").Value = 5896.3335
$ws3.Range("I99").Value = 1866.7142
$ws3.Range("K99").Value = 1866.7142
$ws3.Range("M99").Value = -368.7141999999999

# BSM row 105
$ws3.Range("H105").Value = 2734.2307
$ws3.Range("I105").Value = 2616.111
$ws3.Range("K105").Value = 2616.111
$ws3.Range("M105").Value = -869.1109999999999

# BSM row 134
$ws3.Range("H134").Value = 1761.1356
$ws3.Range("I134").Value = 1752.8909
$ws3.Range("J134").Value = 1874.5
$ws3.Range("K134").Value = 5258.6727
$ws3.Range("L134").Value = 5623.5
$ws3.Range("M134").Value = -2723.6727
$ws3.Range("N134").Value = -10693.5

# CRP row 22
$ws4.Range("H22").Value = 599
$ws4.Range("I22").Value = 599
$ws4.Range("K22").Value = 599
$ws4.Range("M22").Value = -249

# CRP row 31
$ws4.Range("H31").Value = 19643.982
$ws4.Range("I31").Value = 25911.834
$ws4.Range("J31").Value = 4158.706
$ws4.Range("K31").Value = 25911.834
$ws4.Range("L31").Value = 4158.706
$ws4.Range("M31").Value = -25616.834
$ws4.Range("N31").Value = -4748.706

# CRP row 34
$ws4.Range("H34").Value = 19643.982
$ws4.Range("I34").Value = 25911.834
$ws4.Range("J34").Value = 4158.706
$ws4.Range("K34").Value = 25911.834
$ws4.Range("L34").Value = 4158.706
$ws4.Range("M34").Value = -25709.834
$ws4.Range("N34").Value = -4562.706

# CRP row 86
$ws4.Range("H86").Value = 12399
$ws4.Range("I86").Value = 14000
$ws4.Range("J86").Value = 11998.75
$ws4.Range("K86").Value = 14000
$ws4.Range("L86").Value = 11998.75
$ws4.Range("M86").Value = -12877
$ws4.Range("N86").Value = -14244.75

# CRP row 89
$ws4.Range("H89").Value = 12399
$ws4.Range("I89").Value = 14000
$ws4.Range("J89").Value = 11998.75
$ws4.Range("K89").Value = 70000
$ws4.Range("L89").Value = 59993.75
$ws4.Range("M89").Value = -64384
$ws4.Range("N89").Value = -71225.75

# CRP row 132
$ws4.Range("H132").Value = 134199.19
$ws4.Range("I132").Value = 180520.42
$ws4.Range("K132").Value = 541561.26
$ws4.Range("M132").Value = -539031.26

# CRP row 134
$ws4.Range("H134").Value = 24209.564
$ws4.Range("I134").Value = 21450.967
$ws4.Range("J134").Value = 34899.125
$ws4.Range("K134").Value = 64352.901
$ws4.Range("L134").Value = 104697.375
$ws4.Range("M134").Value = -61817.901
$ws4.Range("N134").Value = -109767.375

# CUL row 56
$ws5.Range("H56").Value = 14961.333
$ws5.Range("I56").Value = 14961.333
$ws5.Range("K56").Value = 14961.333
$ws5.Range("M56").Value = -14431.333

# CUL row 107
$ws5.Range("H107").Value = 1454.7273
$ws5.Range("J107").Value = 1668.9375
$ws5.Range("L107").Value = 5006.8125
$ws5.Range("N107").Value = -8846.8125

# CUL row 109
$ws5.Range("H109").Value = 1334.125
$ws5.Range("I109").Value = 737.1667
$ws5.Range("K109").Value = 2211.5001
$ws5.Range("M109").Value = -1171.5001

# CUL row 132
$ws5.Range("H132").Value = 1141.25
$ws5.Range("J132").Value = 1500
$ws5.Range("L132").Value = 13500
$ws5.Range("N132").Value = -18560

# CUL row 136
$ws5.Range("H136").Value = 440429.75
$ws5.Range("I136").Value = 626242.8
$ws5.Range("K136").Value = 1878728.4
$ws5.Range("M136").Value = -1873628.4

# CUL row 137
$ws5.Range("H137").Value = 3885.5
$ws5.Range("J137").Value = 4999.75
$ws5.Range("L137").Value = 14999.25
$ws5.Range("N137").Value = -25199.25

# CUL row 138
$ws5.Range("H138").Value = 27787246
$ws5.Range("I138").Value = 35722816
$ws5.Range("K138").Value = 107168448
$ws5.Range("M138").Value = -107163308

# GSM row 62
$ws6.Range("H62").Value = 57500
$ws6.Range("I62").Value = 40000
$ws6.Range("K62").Value = 40000
$ws6.Range("M62").Value = -39314

# GSM row 65
$ws6.Range("H65").Value = 57500
$ws6.Range("I65").Value = 40000
$ws6.Range("K65").Value = 120000
$ws6.Range("M65").Value = -116568

# GSM row 132
$ws6.Range("H132").Value = 2750.375
$ws6.Range("I132").Value = 2748.9211
$ws6.Range("K132").Value = 8246.763300000001
$ws6.Range("M132").Value = -5716.763300000001

# LTW row 7
$ws7.Range("H7").Value = 4375.5
$ws7.Range("I7").Value = 3710.8
$ws7.Range("J7").Value = 7699
$ws7.Range("K7").Value = 3710.8
$ws7.Range("L7").Value = 7699
$ws7.Range("M7").Value = -3598.8
$ws7.Range("N7").Value = -7923

# LTW row 22
$ws7.Range("H22").Value = 1933.5625
$ws7.Range("I22").Value = 617.125
$ws7.Range("J22").Value = 3250
$ws7.Range("K22").Value = 617.125
$ws7.Range("L22").Value = 3250
$ws7.Range("M22").Value = -322.125
$ws7.Range("N22").Value = -3840

# LTW row 27
$ws7.Range("H27").Value = 1933.5625
$ws7.Range("I27").Value = 617.125
$ws7.Range("J27").Value = 3250
$ws7.Range("K27").Value = 617.125
$ws7.Range("L27").Value = 3250
$ws7.Range("M27").Value = -510.125
$ws7.Range("N27").Value = -3464

# LTW row 63
$ws7.Range("H63").Value = 38000
$ws7.Range("J63").Value = 0
$ws7.Range("L63").Value = 0
$ws7.Range("N63").ClearContents()

# LTW row 66
$ws7.Range("H66").Value = 38000
$ws7.Range("J66").Value = 0
$ws7.Range("L66").Value = 0
$ws7.Range("N66").ClearContents()

# LTW row 74
$ws7.Range("H74").Value = 20000
$ws7.Range("J74").Value = 20000
$ws7.Range("L74").Value = 20000
$ws7.Range("N74").Value = -21996

# LTW row 77
$ws7.Range("H77").Value = 20000
$ws7.Range("J77").Value = 20000
$ws7.Range("L77").Value = 60000
$ws7.Range("N77").Value = -69984

# LTW row 126
$ws7.Range("H126").Value = 4375.5
$ws7.Range("I126").Value = 3710.8
$ws7.Range("J126").Value = 7699
$ws7.Range("K126").Value = 11132.4
$ws7.Range("L126").Value = 23097
$ws7.Range("M126").Value = -8662.400000000001
$ws7.Range("N126").Value = -28037

# LTW row 132
$ws7.Range("H132").Value = 21642.27
$ws7.Range("I132").Value = 24379.979
$ws7.Range("K132").Value = 73139.93700000001
$ws7.Range("M132").Value = -70609.93700000001

# WVR row 14
$ws8.Range("H14").Value = 279735.34
$ws8.Range("I14").Value = 429440.2
$ws8.Range("K14").Value = 429440.2
$ws8.Range("M14").Value = -429272.2

# WVR row 75
$ws8.Range("H75").Value = 32750
$ws8.Range("J75").Value = 17500
$ws8.Range("L75").Value = 17500
$ws8.Range("N75").Value = -19372

# WVR row 78
$ws8.Range("H78").Value = 32750
$ws8.Range("J78").Value = 17500
$ws8.Range("L78").Value = 52500
$ws8.Range("N78").Value = -61860

# WVR row 100
$ws8.Range("H100").Value = 1681.5883
$ws8.Range("I100").Value = 994
$ws8.Range("K100").Value = 1988
$ws8.Range("M100").Value = -1447

# WVR row 126
$ws8.Range("H126").Value = 2158.7827
$ws8.Range("I126").Value = 1941.5625
$ws8.Range("J126").Value = 2655.2856
$ws8.Range("K126").Value = 5824.6875
$ws8.Range("L126").Value = 7965.8568
$ws8.Range("M126").Value = -3354.6875
$ws8.Range("N126").Value = -12905.8568

# WVR row 132
$ws8.Range("H132").Value = 1450.08
$ws8.Range("I132").Value = 1448
$ws8.Range("K132").Value = 4344
$ws8.Range("M132").Value = -1814

# WVR row 136
$ws8.Range("H136").Value = 1360.4517
$ws8.Range("I136").Value = 1264.2264
$ws8.Range("J136").Value = 1927.1111
$ws8.Range("K136").Value = 3792.6792
$ws8.Range("L136").Value = 5781.3333
$ws8.Range("M136").Value = -1242.6792
$ws8.Range("N136").Value = -10881.3333

